$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 9 and 10 (swap data) ---
# Row 9
$ws.Range("A9").Value2 = 7
$ws.Range("B9").Value2 = 6814328
$ws.Range("C9").Value2 = "Slovenia Prva Liga"
$ws.Range("D9").Value2 = "Slovenia Prva Liga"
$ws.Range("E9").Value2 = 45137.52083333334
$ws.Range("F9").Value2 = "NK Domzale"
$ws.Range("G9").Value2 = "NK Bravo"
$ws.Range("H9").Value2 = 1
$ws.Range("I9").Value2 = 1
$ws.Range("J9").Value2 = "D"
$ws.Range("K9").Value2 = 2.35
$ws.Range("L9").Value2 = 3.1
$ws.Range("M9").Value2 = 2.9
$ws.Range("N9").Value2 = 2.15
$ws.Range("O9").Value2 = 3.1
$ws.Range("P9").Value2 = 3.3
$ws.Range("Q9").Value2 = -0.25
$ws.Range("R9").Value2 = 1.925
$ws.Range("S9").Value2 = 1.875
$ws.Range("T9").Value2 = 2.25
$ws.Range("U9").Value2 = 1.95
$ws.Range("V9").Value2 = 1.85
$ws.Range("W9").Value2 = -1
$ws.Range("X9").Value2 = 2.1
$ws.Range("Y9").Value2 = -1
$ws.Range("Z9").Value2 = -0.5
$ws.Range("AA9").Value2 = 0.4375
$ws.Range("AB9").Value2 = -0.5
$ws.Range("AC9").Value2 = 0.425

# Row 10
$ws.Range("A10").Value2 = 8
$ws.Range("B10").Value2 = 6814330
$ws.Range("C10").Value2 = "Slovenia Prva Liga"
$ws.Range("D10").Value2 = "Slovenia Prva Liga"
$ws.Range("E10").Value2 = 45137.52083333334
$ws.Range("F10").Value2 = "NK Maribor"
$ws.Range("G10").Value2 = "NK Aluminij"
$ws.Range("H10").Value2 = 1
$ws.Range("I10").Value2 = 0
$ws.Range("J10").Value2 = "H"
$ws.Range("K10").Value2 = 1.363
$ws.Range("L10").Value2 = 4.5
$ws.Range("M10").Value2 = 7
$ws.Range("N10").Value2 = 1.4
$ws.Range("O10").Value2 = 4.5
$ws.Range("P10").Value2 = 7
$ws.Range("Q10").Value2 = -1.25
$ws.Range("R10").Value2 = 1.85
$ws.Range("S10").Value2 = 1.95
$ws.Range("T10").Value2 = 2.75
$ws.Range("U10").Value2 = 1.8
$ws.Range("V10").Value2 = 2
$ws.Range("W10").Value2 = 0.3999999999999999
$ws.Range("X10").Value2 = -1
$ws.Range("Y10").Value2 = -1
$ws.Range("Z10").Value2 = -0.5
$ws.Range("AA10").Value2 = 0.475
$ws.Range("AB10").Value2 = -1
$ws.Range("AC10").Value2 = 1

# --- Update existing rows 138 and 139 (replace with new match data) ---
# Row 138
$ws.Range("A138").Value2 = 136
$ws.Range("B138").Value2 = 8035687
$ws.Range("C138").Value2 = "Slovenia Prva Liga"
$ws.Range("D138").Value2 = "Slovenia Prva Liga"
$ws.Range("E138").Value2 = 45388.41666666666
$ws.Range("F138").Value2 = "NK Rogaska"
$ws.Range("G138").Value2 = "Olimpija Ljubljana"
$ws.Range("H138").Value2 = 2
$ws.Range("I138").Value2 = 3
$ws.Range("J138").Value2 = "A"
$ws.Range("K138").Value2 = 6
$ws.Range("L138").Value2 = 4.333
$ws.Range("M138").Value2 = 1.45
$ws.Range("N138").Value2 = 4.75
$ws.Range("O138").Value2 = 4
$ws.Range("P138").Value2 = 1.6
$ws.Range("Q138").Value2 = 1
$ws.Range("R138").Value2 = 1.8
$ws.Range("S138").Value2 = 2
$ws.Range("T138").Value2 = 3
$ws.Range("U138").Value2 = 2
$ws.Range("V138").Value2 = 1.8
$ws.Range("W138").Value2 = -1
$ws.Range("X138").Value2 = -1
$ws.Range("Y138").Value2 = 0.6000000000000001
$ws.Range("Z138").Value2 = 0
$ws.Range("AA138").Value2 = 0
$ws.Range("AB138").Value2 = 1
$ws.Range("AC138").Value2 = -1

# Row 139
$ws.Range("A139").Value2 = 137
$ws.Range("B139").Value2 = 6814435
$ws.Range("C139").Value2 = "Slovenia Prva Liga"
$ws.Range("D139").Value2 = "Slovenia Prva Liga"
$ws.Range("E139").Value2 = 45388.52083333334
$ws.Range("F139").Value2 = "NK Radomlje"
$ws.Range("G139").Value2 = "FC Koper"
$ws.Range("H139").Value2 = 1
$ws.Range("I139").Value2 = 1
$ws.Range("J139").Value2 = "D"
$ws.Range("K139").Value2 = 2.55
$ws.Range("L139").Value2 = 3.25
$ws.Range("M139").Value2 = 2.55
$ws.Range("N139").Value2 = 2.6
$ws.Range("O139").Value2 = 3.2
$ws.Range("P139").Value2 = 2.55
$ws.Range("Q139").Value2 = 0
$ws.Range("R139").Value2 = 1.9
$ws.Range("S139").Value2 = 1.9
$ws.Range("T139").Value2 = 2.25
$ws.Range("U139").Value2 = 1.9
$ws.Range("V139").Value2 = 1.9
$ws.Range("W139").Value2 = -1
$ws.Range("X139").Value2 = 2.2
$ws.Range("Y139").Value2 = -1
$ws.Range("Z139").Value2 = 0
$ws.Range("AA139").Value2 = 0
$ws.Range("AB139").Value2 = -0.5
$ws.Range("AC139").Value2 = 0.45

# --- Add new rows 140 and 141: full format copy from row 139, then set values ---
$ws.Range("A139:AC139").Copy()
$ws.Range("A140:AC141").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 140
$ws.Range("A140").Value2 = 138
$ws.Range("B140").Value2 = 6837117
$ws.Range("C140").Value2 = "Slovenia Prva Liga"
$ws.Range("D140").Value2 = "Slovenia Prva Liga"
$ws.Range("E140").Value2 = 45388.63541666666
$ws.Range("F140").Value2 = "NS Mura"
$ws.Range("G140").Value2 = "NK Celje"
$ws.Range("H140").Value2 = 1
$ws.Range("I140").Value2 = 3
$ws.Range("J140").Value2 = "A"
$ws.Range("K140").Value2 = 5.25
$ws.Range("L140").Value2 = 4.2
$ws.Range("M140").Value2 = 1.5
$ws.Range("N140").Value2 = 5.75
$ws.Range("O140").Value2 = 4.333
$ws.Range("P140").Value2 = 1.45
$ws.Range("Q140").Value2 = 1.25
$ws.Range("R140").Value2 = 1.8
$ws.Range("S140").Value2 = 2
$ws.Range("T140").Value2 = 2.75
$ws.Range("U140").Value2 = 1.975
$ws.Range("V140").Value2 = 1.825
$ws.Range("W140").Value2 = -1
$ws.Range("X140").Value2 = -1
$ws.Range("Y140").Value2 = 0.45
$ws.Range("Z140").Value2 = -1
$ws.Range("AA140").Value2 = 1
$ws.Range("AB140").Value2 = 0.9750000000000001
$ws.Range("AC140").Value2 = -1

# Row 141
$ws.Range("A141").Value2 = 139
$ws.Range("B141").Value2 = 6814434
$ws.Range("C141").Value2 = "Slovenia Prva Liga"
$ws.Range("D141").Value2 = "Slovenia Prva Liga"
$ws.Range("E141").Value2 = 45389.41666666666
$ws.Range("F141").Value2 = "NK Bravo"
$ws.Range("G141").Value2 = "NK Domzale"
$ws.Range("H141").Value2 = 1
$ws.Range("I141").Value2 = 3
$ws.Range("J141").Value2 = "A"
$ws.Range("K141").Value2 = 1.833
$ws.Range("L141").Value2 = 3.25
$ws.Range("M141").Value2 = 4
$ws.Range("N141").Value2 = 1.75
$ws.Range("O141").Value2 = 3.4
$ws.Range("P141").Value2 = 4.333
$ws.Range("Q141").Value2 = -0.5
$ws.Range("R141").Value2 = 1.775
$ws.Range("S141").Value2 = 2.025
$ws.Range("T141").Value2 = 2.25
$ws.Range("U141").Value2 = 1.775
$ws.Range("V141").Value2 = 2.025
$ws.Range("W141").Value2 = -1
$ws.Range("X141").Value2 = -1
$ws.Range("Y141").Value2 = 3.333
$ws.Range("Z141").Value2 = -1
$ws.Range("AA141").Value2 = 1.025
$ws.Range("AB141").Value2 = 0.7749999999999999
$ws.Range("AC141").Value2 = -1

# --- Add new rows 142-144: partial format copy (A:G and K:AA only; no H,I,J,AB,AC - future fixtures) ---
$ws.Range("A139:G139").Copy()
$ws.Range("A142:G144").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K139:AA139").Copy()
$ws.Range("K142:AA144").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 142
$ws.Range("A142").Value2 = 140
$ws.Range("B142").Value2 = 6919052
$ws.Range("C142").Value2 = "Slovenia Prva Liga"
$ws.Range("D142").Value2 = "Slovenia Prva Liga"
$ws.Range("E142").Value2 = 45392.40625
$ws.Range("F142").Value2 = "NK Celje"
$ws.Range("G142").Value2 = "NK Bravo"
$ws.Range("K142").Value2 = 1.444
$ws.Range("L142").Value2 = 4.2
$ws.Range("M142").Value2 = 6
$ws.Range("N142").Value2 = 1.444
$ws.Range("O142").Value2 = 4.2
$ws.Range("P142").Value2 = 6
$ws.Range("Q142").Value2 = -1.25
$ws.Range("R142").Value2 = 2.025
$ws.Range("S142").Value2 = 1.775
$ws.Range("T142").Value2 = 2.5
$ws.Range("U142").Value2 = 1.875
$ws.Range("V142").Value2 = 1.925
$ws.Range("W142").Value2 = 0
$ws.Range("X142").Value2 = 0
$ws.Range("Y142").Value2 = 0
$ws.Range("Z142").Value2 = 0
$ws.Range("AA142").Value2 = 0

# Row 143
$ws.Range("A143").Value2 = 141
$ws.Range("B143").Value2 = 6860865
$ws.Range("C143").Value2 = "Slovenia Prva Liga"
$ws.Range("D143").Value2 = "Slovenia Prva Liga"
$ws.Range("E143").Value2 = 45392.57291666666
$ws.Range("F143").Value2 = "Olimpija Ljubljana"
$ws.Range("G143").Value2 = "NK Radomlje"
$ws.Range("K143").Value2 = 1.363
$ws.Range("L143").Value2 = 4.5
$ws.Range("M143").Value2 = 7
$ws.Range("N143").Value2 = 1.363
$ws.Range("O143").Value2 = 4.5
$ws.Range("P143").Value2 = 7
$ws.Range("Q143").Value2 = -1.25
$ws.Range("R143").Value2 = 1.825
$ws.Range("S143").Value2 = 1.975
$ws.Range("T143").Value2 = 2.75
$ws.Range("U143").Value2 = 1.9
$ws.Range("V143").Value2 = 1.9
$ws.Range("W143").Value2 = 0
$ws.Range("X143").Value2 = 0
$ws.Range("Y143").Value2 = 0
$ws.Range("Z143").Value2 = 0
$ws.Range("AA143").Value2 = 0

# Row 144
$ws.Range("A144").Value2 = 142
$ws.Range("B144").Value2 = 6847813
$ws.Range("C144").Value2 = "Slovenia Prva Liga"
$ws.Range("D144").Value2 = "Slovenia Prva Liga"
$ws.Range("E144").Value2 = 45393.48958333334
$ws.Range("F144").Value2 = "FC Koper"
$ws.Range("G144").Value2 = "NK Aluminij"
$ws.Range("K144").Value2 = 1.615
$ws.Range("L144").Value2 = 3.8
$ws.Range("M144").Value2 = 4.75
$ws.Range("N144").Value2 = 1.6
$ws.Range("O144").Value2 = 3.8
$ws.Range("P144").Value2 = 4.75
$ws.Range("Q144").Value2 = -0.75
$ws.Range("R144").Value2 = 1.775
$ws.Range("S144").Value2 = 2.025
$ws.Range("T144").Value2 = 2.5
$ws.Range("U144").Value2 = 1.8
$ws.Range("V144").Value2 = 2
$ws.Range("W144").Value2 = 0
$ws.Range("X144").Value2 = 0
$ws.Range("Y144").Value2 = 0
$ws.Range("Z144").Value2 = 0
$ws.Range("AA144").Value2 = 0

Write-Host "Edit complete"